$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Sending cluster): "Resolving-Mac" -> "MuSCs" for all data rows
$ws.Range("A2:A5").Value = "MuSCs"

# Row 2 (originally Resolving-Mac -> ECs)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.06815433333333333
$ws.Range("H2").Value = 0.204463
$ws.Range("M2").Value = 3.390429
$ws.Range("N2").Value = 10.171287
$ws.Range("O2").Value = 0.173121426386348
$ws.Range("P2").Value = 0.173121426386348
$ws.Range("Q2").Value = 0.231072428209
$ws.Range("R2").Value = 2.079651853881
$ws.Range("S2").Value = 0.173121426386348
$ws.Range("T2").Value = 0.173121426386348

# Row 3 (originally Resolving-Mac -> FAPs)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.06815433333333333
$ws.Range("H3").Value = 0.204463
$ws.Range("O3").Value = 0.5936336753560868
$ws.Range("P3").Value = 0.5936336753560868
$ws.Range("Q3").Value = 0.7923477624603333
$ws.Range("R3").Value = 7.131129862143
$ws.Range("S3").Value = 0.5936336753560868
$ws.Range("T3").Value = 0.5936336753560868

# Row 4 (originally Resolving-Mac -> MuSCs)
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.06815433333333333
$ws.Range("H4").Value = 0.204463
$ws.Range("M4").Value = 4.546141666666667
$ws.Range("N4").Value = 13.638425
$ws.Range("O4").Value = 0.2321342018628743
$ws.Range("P4").Value = 0.2321342018628743
$ws.Range("Q4").Value = 0.3098392545305556
$ws.Range("R4").Value = 2.788553290775
$ws.Range("S4").Value = 0.2321342018628743
$ws.Range("T4").Value = 0.2321342018628743

# Row 5 (originally Resolving-Mac -> Resolving-Mac)
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.021752
$ws.Range("H5").Value = 0.06525600000000001
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.021752
$ws.Range("N5").Value = 0.06525600000000001
$ws.Range("O5").Value = 0.001110696394691009
$ws.Range("P5").Value = 0.001110696394691009
$ws.Range("Q5").Value = 0.001482493058666667
$ws.Range("R5").Value = 0.013342437528
$ws.Range("S5").Value = 0.001110696394691009
$ws.Range("T5").Value = 0.001110696394691009
